# Adds two new "batches" of master device rows (10 devices total) for
# mac-address / serial-number test data, mirroring the existing layout
# of xl/worksheets/sheet1.xml (columns: A id, B name, C mac_address,
# D serial_num, F dspec_id, G lang_code, H is_active, I cr_by,
# J cr_dtimes, K eff_dtimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlHAlignLeft = -4131

$rows = @(
    @{ Row = 147; Id = 3000166; Name = "Finger Print Scanner 30"; Mac = "D6-15-AC-80-6B-86"; Serial = "BS563Q2230814"; Dspec = 165 },
    @{ Row = 148; Id = 3000167; Name = "IRIS Scanner 30";         Mac = "6D-58-E2-DF-74-34"; Serial = "BS563Q2230815"; Dspec = 327 },
    @{ Row = 149; Id = 3000168; Name = "Web Camera 30";           Mac = "E2-A8-56-86-15-30"; Serial = "BS563Q2230816"; Dspec = 736 },
    @{ Row = 150; Id = 3000169; Name = "Document Scanner 30";     Mac = "72-E8-B9-FD-63-65"; Serial = "BS563Q2230817"; Dspec = 801 },
    @{ Row = 151; Id = 3000170; Name = "Printer 30";              Mac = "D3-F3-A4-50-AD-12"; Serial = "BS563Q2230818"; Dspec = 920 },
    @{ Row = 152; Id = 3000171; Name = "Finger Print Scanner 31"; Mac = "06-16-D0-0B-A6-E4"; Serial = "BS563Q2230819"; Dspec = 165 },
    @{ Row = 153; Id = 3000172; Name = "IRIS Scanner 31";         Mac = "21-78-45-AC-E9-20"; Serial = "BS563Q2230820"; Dspec = 327 },
    @{ Row = 154; Id = 3000173; Name = "Web Camera 31";           Mac = "3C-E8-87-99-DB-FA"; Serial = "BS563Q2230821"; Dspec = 736 },
    @{ Row = 155; Id = 3000174; Name = "Document Scanner 31";     Mac = "BF-55-53-98-40-08"; Serial = "BS563Q2230822"; Dspec = 801 },
    @{ Row = 156; Id = 3000175; Name = "Printer 31";              Mac = "5A-43-36-46-22-EB"; Serial = "BS563Q2230823"; Dspec = 920 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.Mac
    $ws.Cells.Item($row, 4).Value = $r.Serial
    $ws.Cells.Item($row, 6).Value = $r.Dspec
    $ws.Cells.Item($row, 7).Value = "eng"
    $ws.Cells.Item($row, 8).Value = $true
    $ws.Cells.Item($row, 8).HorizontalAlignment = $xlHAlignLeft
    $ws.Cells.Item($row, 9).Value = "superadmin"
    $ws.Cells.Item($row, 10).Value = "now()"
    $ws.Cells.Item($row, 11).Value = "now()"
}

# Match the saved view state: scrolled down toward the newly-added rows,
# with D145 as the active (last-edited-area) cell.
$excel.ActiveWindow.ScrollRow = 139
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D145").Select()

Write-Output "Added $($rows.Count) new device rows (147-156)"
